$wb = $excel.ActiveWorkbook

# Sheet 1 = "展览" : update "想去人数" (column F) counts
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F4").Value = 92
$ws1.Range("F8").Value = 1168
$ws1.Range("F9").Value = 256
$ws1.Range("F10").Value = 113
$ws1.Range("F11").Value = 10006
$ws1.Range("F14").Value = 5
$ws1.Range("F15").Value = 615
$ws1.Range("F16").Value = 11737
$ws1.Range("F17").Value = 12071
$ws1.Range("F18").Value = 25
$ws1.Range("F19").Value = 93
$ws1.Range("F21").Value = 25

# Sheet 4 = "全部类型" : update "想去人数" (column F) counts
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F4").Value = 92
$ws4.Range("F9").Value = 1168
$ws4.Range("F10").Value = 256
$ws4.Range("F11").Value = 113
$ws4.Range("F12").Value = 10006
$ws4.Range("F15").Value = 5
$ws4.Range("F16").Value = 615
$ws4.Range("F17").Value = 11737
$ws4.Range("F18").Value = 12071
$ws4.Range("F19").Value = 25
$ws4.Range("F20").Value = 93
$ws4.Range("F22").Value = 25
